# Daily attendance processing - 2026-01-22 21:38:18
# Reorders the "Recorded By" (column G) entries so that the literal
# "System" token is moved to the end of its comma-separated list
# (case-sensitive - a lowercase "system" entry is left in place and is
# not treated as the same token). Rows whose list does not contain the
# exact token "System" simply have their two comma-separated entries
# swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cellRef = "G" + $r
    $cell = $ws.Range($cellRef)
    $original = $cell.Text

    if ($original -eq $null -or $original -eq "") {
        continue
    }
    if ($original -notlike "*,*") {
        continue
    }

    $parts = $original -split ", "

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.CompareTo("System") -eq 0) {
            $hasSystem = $true
        }
    }

    $newParts = @()

    if ($hasSystem) {
        foreach ($p in $parts) {
            if ($p.CompareTo("System") -ne 0) {
                $newParts += $p
            }
        }
        $newParts += "System"
    } else {
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $newParts += $parts[$i]
        }
    }

    $newValue = $newParts -join ", "

    if ($newValue -ne $original) {
        $cell.Value = $newValue
    }
}
